# Generate Report for Handoff
# Adds a new localization-status row (163191f6-ff6e-446b-85ce-d28caa2b9388.md)
# to the Overview, zh-cn and de-de sheets/tables.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Overview" (sheet1 / table3) -> new row 3
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.ListObjects.Item(1).ListRows.Add()

$wsOverview.Range("A3").Value = "163191f6-ff6e-446b-85ce-d28caa2b9388.md"
$wsOverview.Range("B3").Value = "e2e\163191f6-ff6e-446b-85ce-d28caa2b9388.md"
$wsOverview.Range("C3").Value = ".md"
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-08-26 12:40:27"
$wsOverview.Range("G3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9d415bcef62a5296eff1ec2a8eba80b52bc9434d/e2e/163191f6-ff6e-446b-85ce-d28caa2b9388.md", "", "", "e2e\163191f6-ff6e-446b-85ce-d28caa2b9388.md")

# ---------------------------------------------------------------------------
# Sheet "zh-cn" (sheet2 / table1) -> new row 3
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.ListObjects.Item(1).ListRows.Add()

$wsZhCn.Range("A3").Value = "163191f6-ff6e-446b-85ce-d28caa2b9388.md"
$wsZhCn.Range("B3").Value = ".md"
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("D3").Value = "e2e"
$wsZhCn.Range("E3").Value = "ht"
$wsZhCn.Range("F3").Value = "'False"
$wsZhCn.Range("G3").Value = "163191f6-ff6e-446b-85ce-d28caa2b9388.428b6423ebf8dd44fb09048845a5bccc6a199735.zh-cn.xlf"
$wsZhCn.Range("H3").Value = "2016-08-26 12:40:22"
$wsZhCn.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Range("K3").Value = "0001-01-01 00:00:00"
$wsZhCn.Range("K3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Range("M3").Value = "'True"
$wsZhCn.Range("O3").Value = "'False"

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9d415bcef62a5296eff1ec2a8eba80b52bc9434d/e2e/163191f6-ff6e-446b-85ce-d28caa2b9388.md", "", "", "163191f6-ff6e-446b-85ce-d28caa2b9388.md")

# ---------------------------------------------------------------------------
# Sheet "de-de" (sheet3 / table2) -> new row 3
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.ListObjects.Item(1).ListRows.Add()

$wsDeDe.Range("A3").Value = "163191f6-ff6e-446b-85ce-d28caa2b9388.md"
$wsDeDe.Range("B3").Value = ".md"
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("D3").Value = "e2e"
$wsDeDe.Range("E3").Value = "ht"
$wsDeDe.Range("F3").Value = "'False"
$wsDeDe.Range("G3").Value = "163191f6-ff6e-446b-85ce-d28caa2b9388.428b6423ebf8dd44fb09048845a5bccc6a199735.de-de.xlf"
$wsDeDe.Range("H3").Value = "2016-08-26 12:40:27"
$wsDeDe.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Range("K3").Value = "0001-01-01 00:00:00"
$wsDeDe.Range("K3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Range("M3").Value = "'True"
$wsDeDe.Range("O3").Value = "'False"

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9d415bcef62a5296eff1ec2a8eba80b52bc9434d/e2e/163191f6-ff6e-446b-85ce-d28caa2b9388.md", "", "", "163191f6-ff6e-446b-85ce-d28caa2b9388.md")
